$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.232.77"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.865.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.95%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.02"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2832"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06508"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.37"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.42%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07854"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "97.36"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.869.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.092"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6718"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "279.94"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.229.22"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.30%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.0000"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.486"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.30%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.64"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.118.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.57%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.000007269"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.144"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.63%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.72"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.02%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.925"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.70%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.376"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09644"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.394"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.36%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.093"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04698"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.115"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7048"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.63%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01848"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.534"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.236"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -6.65%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.05"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.943"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -2.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8448"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.38%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4166"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.95%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.82"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.179"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.222"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "934.99"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -5.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.07"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.06%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.56%  "

